# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

# Row -> new value map (column F) to apply to both sheets.
$updates = @{
    2  = 1088
    3  = 796
    4  = 268
    5  = 42
    6  = 1104
    8  = 1929
    9  = 7151
    10 = 897
    11 = 406
    12 = 340
    13 = 124
    14 = 394
    15 = 149
    16 = 7084
    17 = 296
    18 = 1328
    19 = 147
    22 = 135
    23 = 294
    24 = 132
    26 = 15
    32 = 89
    33 = 54
    34 = 72
    35 = 35
    36 = 72
    37 = 65
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
